$wb = $excel.ActiveWorkbook

# Sheets: 1 = protocole, 2 = Feuil1
$wsProtocole = $wb.Worksheets.Item(1)
$wsFeuil1 = $wb.Worksheets.Item(2)

# Update the two accented labels to their non-accented equivalents
# (this also causes the shared-string table to be rebuilt / reordered
# to match the target workbook)
$wsFeuil1.Range("A1").Value = "Chaine de caractere"
$wsFeuil1.Range("A4").Value = "Booleen"

# Feuil1 becomes the active sheet / tab, with A2 selected
$wsFeuil1.Activate()
$wsFeuil1.Range("A2").Select()
